$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.503.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.924.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.95"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4846"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.21%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08190"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.027"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.86"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.922.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.058"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.233"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.56"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.08%  "

$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06783"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.523.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.639"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.80"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.182"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.148.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.708"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +10.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.90"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.13"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.78%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.031"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09586"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.561"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.27%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.570"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.396"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.186"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.32%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.049"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.83"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1868"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.57%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.281"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07631"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5603"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.01"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.441"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.30%  "
